# Tidsregistrering i PTE projektet - Nada Omer
# Add three new time-registration rows (58-60) for "Implementer" work on
# GUI programmering / mellemregninger-logik, 22-24 March 2017.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 58: 22-03-2017, GUI programmering ---------------------------------
$ws.Range("A58").Value = 42816
$ws.Range("B58").Value = "NO"
$ws.Range("C58").Value = "Nada H. A. Omer"
$ws.Range("E58").Value = "Implementer"
$ws.Range("F58").Value = "GUI programmering"
$ws.Range("G58").Value = 0.33680555555555558
$ws.Range("H58").Value = 0.70486111111111116
$ws.Range("I58").Value = "8 time : 50 min."

# --- Row 59: 23-03-2017, GUI programmering + Logic for mellemregninger -----
$ws.Range("A59").Value = 42817
$ws.Range("B59").Value = "NO"
$ws.Range("C59").Value = "Nada H. A. Omer"
$ws.Range("E59").Value = "Implementer"
$ws.Range("F59").Value = "GUI programmering + Logic for mellemregniner formler"
$ws.Range("G59").Value = 0.33680555555555558
$ws.Range("H59").Value = 0.63541666666666663
$ws.Range("I59").Value = "7 time : 10 min."

# --- Row 60: 24-03-2017, Logic for mellemregninger --------------------------
$ws.Range("A60").Value = 42818
$ws.Range("B60").Value = "NO"
$ws.Range("C60").Value = "Nada H. A. Omer"
$ws.Range("E60").Value = "Implementer"
$ws.Range("F60").Value = " Logic for mellemregniner formler"
$ws.Range("G60").Value = 0.33680555555555558
$ws.Range("H60").Value = 0.54166666666666663
$ws.Range("I60").Value = "5 time : 00 min."

# Carry the formatting (date/time number formats, vertical-centred bold
# styles, etc.) from the previous last row (57) down onto the new rows,
# without pulling in the unused column D.
$ws.Range("A57:C57").Copy()
$ws.Range("A58:C60").PasteSpecial(-4122)
$ws.Range("E57:I57").Copy()
$ws.Range("E58:I60").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The "Deltagere" validation list on column C covered C3:C57; extend it to
# cover the new rows (C3:C60). Recreate it (Delete+Add) so the sqref
# actually grows, and re-add the "Gyldige roller" validation on column E
# afterwards so the two validations keep their original relative order.
$ws.Range("C3:C57").Validation.Delete()
$ws.Range("E3:E126").Validation.Delete()
$ws.Range("C3:C60").Validation.Add(3, 1, 1, "=Deltagere")
$ws.Range("E3:E126").Validation.Add(3, 1, 1, "=GyldigeRoller")

# Leave the selection on F52, matching where the author ended up editing.
$ws.Range("F52").Select() | Out-Null
